# 🚌 141: 31/12 15:43 LP1912+6203+6173
# Append freshly-scraped rows to the three schedule sheets and refresh the
# "Última actualización" / "Total filas" header cells on each of them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912": columns A Hora_Scrap? -> actually A=(blank) B=Hora_Scrap
#                 C=Hora_Llegada D=Línea E=Minutos F=Parada G=Fecha
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2, 1).Value = "Última actualización: 31/12/2025 12:43:43"
$ws1.Cells.Item(3, 1).Value = "Total filas: 946"

$sheet1Rows = @(
    @("12:43:32", "12:51", "15_ABASTO", 8, "LP1912", "31/12/2025"),
    @("12:43:32", "12:58", "16_SANTA ANA", 15, "LP1912", "31/12/2025"),
    @("12:43:32", "13:01", "215C_EL PATO", 18, "LP1912", "31/12/2025"),
    @("12:43:32", "13:03", "23_HERNANDEZ", 20, "LP1912", "31/12/2025"),
    @("12:43:32", "13:06", "14_ABASTO", 23, "LP1912", "31/12/2025"),
    @("12:43:32", "13:11", "16_SANTA ANA", 28, "LP1912", "31/12/2025"),
    @("12:43:32", "13:18", "11_ETCHEVERRY", 35, "LP1912", "31/12/2025"),
    @("12:43:32", "13:21", "16_SANTA ANA", 38, "LP1912", "31/12/2025"),
    @("12:43:32", "13:21", "17_ROMERO", 38, "LP1912", "31/12/2025"),
    @("12:43:32", "13:30", "10_OLMOS", 47, "LP1912", "31/12/2025"),
    @("12:43:32", "13:31", "16_P MOR-SANTA ANA", 48, "LP1912", "31/12/2025"),
    @("12:43:32", "13:33", "23_HERNANDEZ", 50, "LP1912", "31/12/2025"),
    @("12:43:32", "13:51", "15_ABASTO", 68, "LP1912", "31/12/2025"),
    @("12:43:32", "14:04", "23_HERNANDEZ", 81, "LP1912", "31/12/2025")
)

$startRow1 = 934
for ($i = 0; $i -lt $sheet1Rows.Count; $i++) {
    $r = $startRow1 + $i
    $row = $sheet1Rows[$i]
    $ws1.Cells.Item($r, 2).Value = $row[0]
    $ws1.Cells.Item($r, 3).Value = $row[1]
    $ws1.Cells.Item($r, 4).Value = $row[2]
    $ws1.Cells.Item($r, 5).Value = $row[3]
    $ws1.Cells.Item($r, 6).Value = $row[4]
    $ws1.Cells.Item($r, 7).Value = $row[5]
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215": columns A=(blank) B=Fecha C=Hora_Scrap D=Hora_Llegada
#                     E=Línea F=Minutos G=Parada
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2, 1).Value = "Última actualización: 31/12/2025 12:43:43"
$ws2.Cells.Item(3, 1).Value = "Total filas: 72"

$ws2.Cells.Item(73, 2).Value = "31/12/2025"
$ws2.Cells.Item(73, 3).Value = "12:43:32"
$ws2.Cells.Item(73, 4).Value = "13:01"
$ws2.Cells.Item(73, 5).Value = "215C_EL PATO"
$ws2.Cells.Item(73, 6).Value = 18
$ws2.Cells.Item(73, 7).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "6203-6173": columns A=(blank) B=Fecha C=Hora_Scrap D=Hora_Llegada
#                    E=Línea F=Minutos G=Parada
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2, 1).Value = "Última actualización: 31/12/2025 12:43:43"
$ws3.Cells.Item(3, 1).Value = "Total filas: 113"

$sheet3Rows = @(
    @("31/12/2025", "12:43:43", "13:08", "215B_LP-P MOR-1 Y 57", 25, "L6173"),
    @("31/12/2025", "12:43:43", "13:13", "215A_LA PLATA", 30, "L6173"),
    @("31/12/2025", "12:43:38", "13:53", "215C_LA PLATA", 70, "L6203")
)

$startRow3 = 112
for ($i = 0; $i -lt $sheet3Rows.Count; $i++) {
    $r = $startRow3 + $i
    $row = $sheet3Rows[$i]
    $ws3.Cells.Item($r, 2).Value = $row[0]
    $ws3.Cells.Item($r, 3).Value = $row[1]
    $ws3.Cells.Item($r, 4).Value = $row[2]
    $ws3.Cells.Item($r, 5).Value = $row[3]
    $ws3.Cells.Item($r, 6).Value = $row[4]
    $ws3.Cells.Item($r, 7).Value = $row[5]
}
